# edit.ps1
# Applies the "Add files via upload" revision to TCOtoliths.xlsx:
#  - Rows that previously had "Mounted? = Y" and "Sampled? = Y" (columns I/J)
#    are reworked: those two cells are cleared and a note "RESAMPLE" is
#    written into column L instead.
#  - Row 319 (previously "Argyropelecus olfersii" / TCAO016) is reclassified
#    as "Argyropelecus aculeatus" / TCAA001, with a note (partly italic)
#    explaining the original identification.
#  - The view is scrolled/zoomed to where the author left off, with
#    L320 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the "Mounted?"/"Sampled?" Y/Y pair (columns I & J) is replaced
# by a "RESAMPLE" note in column L.
$resampleRows = @(7,160,161,162,164,167,169,170,172,174,200,202,203,205,207,212,213,215,216,219,291)

foreach ($r in $resampleRows) {
    $ws.Cells.Item($r, 9).ClearContents()   # column I - Mounted?
    $ws.Cells.Item($r, 10).ClearContents()  # column J - Sampled?
    $ws.Cells.Item($r, 12).Value = "RESAMPLE"  # column L - Note:
}

# Row 204 already carried a column L note ("Used as voucher"); only the
# Mounted?/Sampled? Y/Y pair needs to be cleared there, the note is kept.
$ws.Cells.Item(204, 9).ClearContents()
$ws.Cells.Item(204, 10).ClearContents()

# Row 319: specimen was originally logged as Argyropelecus olfersii
# (code TCAO016) and has been reclassified as Argyropelecus aculeatus
# (code TCAA001), with an explanatory note.
$ws.Range("B319").Value = "Argyropelecus aculeatus"
$ws.Range("H319").Value = "TCAA001"

$noteCell = $ws.Range("L319")
$noteCell.Value = "Originally identified as A. olfersii"
$italicRun = $noteCell.Characters(26, 11)
$italicRun.Font.Italic = $true

# Restore the view state (scroll position, zoom, active selection) to
# match where the author left the workbook.
$excel.ActiveWindow.Zoom = 85
$ws.Range("L320").Select()
